$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5516
$ws1.Range("F7").Value = 5
$ws1.Range("F10").Value = 2459
$ws1.Range("F12").Value = 82
$ws1.Range("F14").Value = 2304
$ws1.Range("F15").Value = 222

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 5516
$ws4.Range("F8").Value = 5
$ws4.Range("F12").Value = 2459
$ws4.Range("F14").Value = 82
$ws4.Range("G14").Value = 45
$ws4.Range("F17").Value = 2304
$ws4.Range("F18").Value = 222
